$d = $word.ActiveDocument
$ns = ' xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

$p0 = $d.Paragraphs(1)
$xml0 = '<w:p__NS__ w14:paraId="583C0996" w14:textId="445027AB" w:rsidR="00C844BC" w:rsidRDefault="00282E97"><w:pPr><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>Journalistisch en redactionele vaardigheden</w:t></w:r></w:p>'
$xml0 = $xml0.Replace("__NS__", $ns)
$p0.Range.InsertXML($xml0)

$p1 = $d.Paragraphs(2)
$xml1 = '<w:p__NS__ w14:paraId="2E9FF15D" w14:textId="0C39E4C7" w:rsidR="00282E97" w:rsidRDefault="00282E97"><w:pPr><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>Jort S</w:t></w:r><w:r w:rsidR="00D844F5"><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>i</w:t></w:r><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>emes</w:t></w:r></w:p>'
$xml1 = $xml1.Replace("__NS__", $ns)
$p1.Range.InsertXML($xml1)

$p2 = $d.Paragraphs(3)
$xml2 = '<w:p__NS__ w14:paraId="269EB761" w14:textId="77777777" w:rsidR="00D844F5" w:rsidRDefault="00D844F5"><w:pPr><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr></w:p>'
$xml2 = $xml2.Replace("__NS__", $ns)
$p2.Range.InsertXML($xml2)

$p3 = $d.Paragraphs(4)
$xml3 = '<w:p__NS__ w14:paraId="0341E97B" w14:textId="77777777" w:rsidR="00D844F5" w:rsidRDefault="00D844F5"><w:pPr><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr></w:p>'
$xml3 = $xml3.Replace("__NS__", $ns)
$p3.Range.InsertXML($xml3)

$p6 = $d.Paragraphs(7)
$xml6 = '<w:p__NS__ w14:paraId="4E325F75" w14:textId="6C595000" w:rsidR="00D844F5" w:rsidRDefault="00D844F5"><w:pPr><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve">Intro </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>icke</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$xml6 = $xml6.Replace("__NS__", $ns)
$p6.Range.InsertXML($xml6)

$p9 = $d.Paragraphs(10)
$xml9 = '<w:p__NS__ w14:paraId="21D5A89A" w14:textId="7703ED91" w:rsidR="00D844F5" w:rsidRDefault="00D844F5"><w:pPr><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>Icke’s</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve"> Ruimteschip:</w:t></w:r></w:p>'
$xml9 = $xml9.Replace("__NS__", $ns)
$p9.Range.InsertXML($xml9)

$p10 = $d.Paragraphs(11)
$xml10 = '<w:p__NS__ w14:paraId="3F3541C3" w14:textId="29D30208" w:rsidR="00D844F5" w:rsidRPr="00D844F5" w:rsidRDefault="00D844F5"><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r w:rsidRPr="00D844F5"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve">De ruimteschepen die wij bouwen op aarde zijn extreem primitief, maar die zien er dus uit zo International Space Station. Dat is gewoon een </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="nl-NL"/></w:rPr><w:t>een</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve"> groot bierblikje met een stuk of wat Mensen erin.”</w:t></w:r></w:p>'
$xml10 = $xml10.Replace("__NS__", $ns)
$p10.Range.InsertXML($xml10)

$p11 = $d.Paragraphs(12)
$xml11 = '<w:p__NS__ w14:paraId="0FD66A9C" w14:textId="3522377F" w:rsidR="00D844F5" w:rsidRDefault="00D844F5"><w:pPr><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="nl-NL"/></w:rPr><w:t>“</w:t></w:r><w:r w:rsidRPr="00D844F5"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve">want voor Mensen is de ruimte het gevaarlijkste. Wat? Willekeurig welke plek op onze planeet of het nou de topje van Mount Everest is of op de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="nl-NL"/></w:rPr><w:t>zuidste</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve"> puntje van de Zuidpool is altijd nog prettiger en behaaglijker dan willekeurig. Welk.in de ruimte? Als je mij op het zuidelijkste puntje van de Zuidpool neerzet en Ik heb een behoorlijke winterjas aan, dan kan ik het Misschien nog wel een paar dagen uithouden Als je mij In de ruimte zet. Die witte jas of geen winterjas, dan hou ik het geen 3 minuten uit en. Dan ben ik morsdood</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="nl-NL"/></w:rPr><w:t>”</w:t></w:r></w:p>'
$xml11 = $xml11.Replace("__NS__", $ns)
$p11.Range.InsertXML($xml11)

$p13 = $d.Paragraphs(14)
$xml13 = '<w:p__NS__ w14:paraId="1F25EA85" w14:textId="5B9CB8ED" w:rsidR="00E45734" w:rsidRDefault="00E45734" w:rsidP="00E45734"><w:pPr><w:pBdr><w:bottom w:val="single" w:sz="6" w:space="1" w:color="auto"/></w:pBdr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r w:rsidRPr="00122E54"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="nl-NL"/></w:rPr><w:t>het zou inderdaad kunnen zijn dat die zwerm waar ik het over heb, de massa verdeling in die zwerm zodanig is dat het verandering van die van die massa beding het ruimteschip Voortstuwt</w:t></w:r><w:r w:rsidRPr="00D574EC"><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>.</w:t></w:r><w:r w:rsidRPr="00E45734"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidRPr="00782222"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve">Maar Misschien is het wel heel wat anders. Ik bedoel, Wij zijn tenslotte een </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="nl-NL"/></w:rPr><w:t>een</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve"> geëvolueerde apensoort en Misschien is ons brein zijn onze hersens zelf helemaal niet in staat voldoende te begrijpen van het heelal onduidelijk gebruik van te maken dat het. Ik heb geen idee.</w:t></w:r></w:p>'
$xml13 = $xml13.Replace("__NS__", $ns)
$p13.Range.InsertXML($xml13)
